$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New date column CK is column 89 (one past CJ = column 88), date serial 45994 (2025-12-03)
$dateCol = 89
$prevCol = 88
$xlPasteFormats = -4122

# Row 1: header date cell - set the new date value, then copy formatting from CJ1
# (value must be set first so the per-player formula ranges pick up the new
# column as "non-blank" when the workbook recalculates)
$ws.Cells.Item(1, $dateCol).Value = 45994
$ws.Cells.Item(1, $prevCol).Copy() | Out-Null
$ws.Cells.Item(1, $dateCol).PasteSpecial($xlPasteFormats)

# Attendance values for the new date (column CK) per player row.
# Row 12 is intentionally omitted: that player's data doesn't extend this far (row ends at column AX).
$values = @{
    2  = "P"
    3  = "P"
    4  = "P"
    5  = "P"
    6  = "RH"
    7  = "P"
    8  = "P"
    9  = "P"
    10 = "P"
    11 = "P"
    13 = "B"
    14 = "P"
    15 = "P"
    16 = "P"
    17 = "P"
    18 = "P"
    19 = "P"
    20 = "P"
    22 = "P"
    23 = "A"
    24 = "P"
    25 = "B"
    26 = "M"
    27 = "P"
    28 = "P"
    29 = "P"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, $dateCol).Value = $values[$row]
    $ws.Cells.Item($row, $prevCol).Copy() | Out-Null
    $ws.Cells.Item($row, $dateCol).PasteSpecial($xlPasteFormats)
}

# Row 21 gets a blank (but styled) cell at CK21, matching the rest of that trailing blank run
$ws.Cells.Item(21, $dateCol).Value = ""
$ws.Cells.Item(21, $prevCol).Copy() | Out-Null
$ws.Cells.Item(21, $dateCol).PasteSpecial($xlPasteFormats)

# Update selection to match the new active cell
$ws.Range("CM25").Select()
